# Weekly price update: a new "Poroto granado" price record for
# Macroferia Regional de Talca is inserted as row 104, pushing the
# existing rows 104-177 down to 105-178 (dimension grows to A1:R178).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 104, shifting rows 104:177 -> 105:178.
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new weekly record.
$ws.Range("A104").Value = 5
$ws.Range("B104").Value = "Macroferia Regional de Talca"
$ws.Range("C104").Value = "Maule"
$ws.Range("D104").Value = 44957
$ws.Range("E104").Value = 7
$ws.Range("F104").Value = 100112030
$ws.Range("G104").Value = "Poroto granado"
$ws.Range("H104").Value = "Sin especificar"
$ws.Range("I104").Value = "Primera"
$ws.Range("J104").Value = 300
$ws.Range("K104").Value = 34000
$ws.Range("L104").Value = 34000
$ws.Range("M104").Value = 34000
$ws.Range("N104").Value = "$/saco 25 kilos"
$ws.Range("O104").Value = "Región del Maule"
$ws.Range("P104").Value = 1360
$ws.Range("Q104").Value = 25
$ws.Range("R104").Value = "Hortaliza"
